$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns C and D are no longer part of the spline coefficient table
# (it shrank from 4 columns to 2) - drop them entirely.
$ws.Range("C1:D3").Clear()

# Recomputed spline coefficients (columns A and B only). Values are stored
# as literal text (same as the rest of this generated table), so prefix
# each numeric-looking literal with an apostrophe to force text storage.
$ws.Range("A1").Value = "'-0.333333333333333"
$ws.Range("B1").Value = "'0.666666666666667"

$ws.Range("A2").Value = "'-0.333333333333333"
$ws.Range("B2").Value = "'0.666666666666667"

$ws.Range("A3").Value = "'6.34413156928661e-17"
$ws.Range("B3").Value = "'3.0"
